# Adds a new "google_search" intent column (J) with its example phrases,
# plus a couple of extra phrases for the existing manage_alarm / manage_timer
# intents (cancel my alarm / cancel timer / cancel my timer).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New example phrases for existing intents (manage_alarm / manage_timer) ---
$ws.Range("I11").Value = "cancel timer"
$ws.Range("I12").Value = "cancel my timer"
$ws.Range("H10").Value = "cancel my alarm"

# --- New intent column: google_search ---
$ws.Range("J1").Value = "google_search"
$ws.Range("J1").HorizontalAlignment = -4108  # xlCenter, matches header style of A1:I1

$ws.Range("J2").Value = "google search"
$ws.Range("J3").Value = "search on google"
$ws.Range("J4").Value = "look up"

# Leave the active cell/selection where the author left it while working
$ws.Range("J5").Select()
